$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.734.90'
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").Value = '3.748.05'
$ws.Range("E3").Value = '  +1.94%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.58'
$ws.Range("E5").Value = '  +1.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.33'
$ws.Range("E6").Value = '  +1.67%  '

$ws.Range("D7").Value = '3.745.72'
$ws.Range("E7").Value = '  +1.92%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  +2.66%  '

$ws.Range("E10").Value = '  +4.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.33'
$ws.Range("E11").Value = '  +3.11%  '

$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.29'
$ws.Range("E13").Value = '  +1.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").Value = '  +3.51%  '

$ws.Range("D15").Value = '4.375.31'

$ws.Range("D16").Value = '3.753.78'
$ws.Range("E16").Value = '  +1.91%  '

$ws.Range("D17").Value = '68.785.75'
$ws.Range("E17").Value = '  +2.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.28'
$ws.Range("E18").Value = '  +2.80%  '

$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  +1.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.81'
$ws.Range("E21").Value = '  +19.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '495.33'
$ws.Range("E22").Value = '  +2.45%  '

$ws.Range("E23").Value = '  +1.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000153'
$ws.Range("E24").Value = '  +9.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.29'
$ws.Range("E25").Value = '  +0.71%  '

$ws.Range("E26").Value = '  +1.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.34'
$ws.Range("E27").Value = '  +1.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.32'
$ws.Range("E28").Value = '  +3.67%  '

$ws.Range("E29").Value = '  +0.44%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.99'
$ws.Range("E30").Value = '  +2.88%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.51'
$ws.Range("E31").Value = '  +6.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.93'
$ws.Range("E32").Value = '  +2.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.86'
$ws.Range("E33").Value = '  +0.72%  '

$ws.Range("D34").Value = '3.895.20'
$ws.Range("E34").Value = '  +2.02%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.109'
$ws.Range("E35").Value = '  +2.11%  '

$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.682.22'
$ws.Range("E36").Value = '  +1.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  +2.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.85'
$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("E40").Value = '  +1.15%  '

$ws.Range("E41").Value = '  +0.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '437.62'
$ws.Range("E42").Value = '  -0.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '48.84'
$ws.Range("E43").Value = '  +0.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.91'
$ws.Range("E44").Value = '  +5.64%  '

$ws.Range("E45").Value = '  +2.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.47'
$ws.Range("E46").Value = '  +2.24%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.72'
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.56'
$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("D50").Value = '2.787.22'
$ws.Range("E50").Value = '  +1.34%  '

$ws.Range("E51").Value = '  +2.88%  '
